# Auto push commit: duplicate the "4b" sheet into a new "temp_2" sheet
# (placed right after "temp"), then revise the live transaction numbers
# on the original "4b" sheet.

$wb = $excel.ActiveWorkbook

# --- 1. Duplicate "4b" -> new sheet inserted immediately before it,
#        then rename the duplicate to "temp_2" so the final tab order
#        becomes: temp, temp_2, 4b, adjusted entry, 2015 -----------------
$source = $wb.Worksheets.Item("4b")
$source.Copy($source)
$duplicate = $wb.Worksheets.Item(2)
$duplicate.Name = "temp_2"

# --- 2. Update the figures on the original "4b" sheet ------------------
$ws = $wb.Worksheets.Item("4b")
$ws.Activate()
$ws.Range("H5").Select()

$ws.Range("B4").Value = 57000
$ws.Range("E4").Value = 32500
$ws.Range("H4").Formula = "=B4+E4"

$ws.Range("D5").ClearContents()
$ws.Range("G5").ClearContents()
$ws.Range("B6").ClearContents()
$ws.Range("C6").ClearContents()
$ws.Range("I6").ClearContents()
$ws.Range("G7").ClearContents()
$ws.Range("J7").ClearContents()
$ws.Range("B8").ClearContents()
$ws.Range("C8").ClearContents()

# --- 3. Leave "temp_2" as the workbook's active tab ---------------------
$wb.Worksheets.Item("temp_2").Activate()
